$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$lo = $ws.ListObjects.Item(1)

# Insert a new worksheet column before column F (6), shifting _price/_spritePath
# one column to the right. This is equivalent to selecting the "_stack" column
# in the table and choosing "Insert > Table Column to the Right".
$ws.Columns.Item(6).Insert()

# Grow the table so it now spans the new column too.
$lo.Resize($ws.Range("A1:H19"))

# Re-assert the header text for every column from F onward so the ListObject
# picks up the right names in the right slots (F is the brand new column).
$ws.Range("F1").Value = "_maxstack"
$ws.Range("G1").Value = "_price"
$ws.Range("H1").Value = "_spritePath"

# Column widths for the new layout (F/G/H), matching the widths Excel would
# compute after the insert + user resize.
$ws.Columns.Item(6).ColumnWidth = 12.75
$ws.Columns.Item(7).ColumnWidth = 19.125
$ws.Columns.Item(8).ColumnWidth = 17.5

# _stack (E) now holds a flat "1" for every item, while the old _stack values
# (all 20) move over to the new _maxstack column (F).
for ($r = 2; $r -le 19; $r++) {
    $ws.Range("E" + $r).Value = 1
    $ws.Range("F" + $r).Value = 20
}

# The _price values that used to live in F14:F19 now belong in G14:G19 (they
# already shifted there physically via the column insert above, but re-assert
# them explicitly for safety/clarity).
$prices = @{14=100; 15=200; 16=300; 17=400; 18=500; 19=600}
foreach ($r in $prices.Keys) {
    $ws.Range("G" + $r).Value = $prices[$r]
}

# Update the active selection to match the authored workbook.
$ws.Range("F9").Select()
